$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("étape 2")
$ws.Activate()

# --- Row 8 ---
$ws.Range("H8").Value = 1

# --- Row 9 : taller row + new values + new comments ---
$ws.Rows(9).RowHeight = 60
$ws.Range("H9").Value = 5
$ws.Range("J9").Value = "largeur des cartes, div image qui sortait de la largeur de la carte"
$ws.Range("K9").Value = "largeur des cartes fixée à 30% de l'espace hébergement, largeur div image en mode auto"

# --- Row 10 ---
$ws.Range("H10").Value = 1

# --- Row 11 : was a blank spacer row, now becomes a real data row ---
$ws.Rows(11).RowHeight = 15.75
$ws.Range("A11").Value = "réaliser une card css"
$ws.Range("B11").Value = "autres"
$ws.Range("C11").Value = "séance guidée"
$ws.Range("D11").Value = 44261
$ws.Range("E11").Value = 44261
$ws.Range("F11").Value = 1
$ws.Range("H11").Value = 1

# Re-assert the G/I shared formulas across their final live range (7:11) so the
# shared-formula grouping survives the upcoming row deletion below.
$ws.Range("G7:G11").Formula = "=F7*30"
$ws.Range("I7:I11").Formula = "=H7*30"

# --- Remove the now-unused rows 12:16 (old blank rows + old totals row) ---
$ws.Rows("12:16").Delete()

# --- Selection moves to A14 ---
$ws.Range("A14").Select()
